# Apply the NATMI LR-pair results update ("Natmi following Dr Hou advice"):
# expands the sending/target cluster cross join from 12 to 18 rows
# (M1 / M2 / sCs targets for every sending cluster) and refreshes all
# derived statistic columns (E:T) for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol, D Target cluster,
# E..T numeric NATMI statistics
$data = @(
    @("ECs", "Col9a3", "Mag", "M1", 2, 1, 1.235259, 2.470518, 0.2325174299743795, 0.1732109743374089, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.34147499796, 2.04884998776, 0.04511863359246258, 0.03825794222969786),
    @("ECs", "Col9a3", "Mag", "M2", 2, 1, 1.235259, 2.470518, 0.2325174299743795, 0.1732109743374089, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 0.77699849865, 4.661990991900001, 0.102663769739855, 0.0870528262727914),
    @("ECs", "Col9a3", "Mag", "sCs", 2, 1, 1.235259, 2.470518, 0.2325174299743795, 0.1732109743374089, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.6413069445120001, 2.565227778048, 0.08473502664206199, 0.04790020583491966),
    @("FAPs", "Col9a3", "Mag", "M1", 3, 1, 0.5501983333333333, 1.650595, 0.103565893827006, 0.1157251913106706, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.1520968272666667, 1.3688714454, 0.02009634983825275, 0.02556078043334561),
    @("FAPs", "Col9a3", "Mag", "M2", 3, 1, 0.5501983333333333, 1.650595, 0.103565893827006, 0.1157251913106706, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 0.3460839216388889, 3.11475529475, 0.04572760449799216, 0.0581614705020316),
    @("FAPs", "Col9a3", "Mag", "sCs", 3, 1, 0.5501983333333333, 1.650595, 0.103565893827006, 0.1157251913106706, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.28564536832, 1.71387220992, 0.03774193949076114, 0.03200294037529344),
    @("M1", "Col9a3", "Mag", "M1", 3, 1, 0.7294966666666666, 2.18849, 0.1373158909250691, 0.1534376536530703, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.2016620585333333, 1.8149585268, 0.0266453373828939, 0.03389051364542637),
    @("M1", "Col9a3", "Mag", "M2", 3, 1, 0.7294966666666666, 2.18849, 0.1373158909250691, 0.1534376536530703, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 0.458865561611111, 4.1297900545, 0.06062929135724442, 0.07711509884556243),
    @("M1", "Col9a3", "Mag", "sCs", 3, 1, 0.7294966666666666, 2.18849, 0.1373158909250691, 0.1534376536530703, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.37873132544, 2.27238795264, 0.05004126218493078, 0.04243204116208151),
    @("M2", "Col9a3", "Mag", "M1", 3, 1, 1.715756666666667, 5.14727, 0.3229633061525894, 0.3608812612892173, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.4743037729333333, 4.268733956399999, 0.06266912151796367, 0.079709582484587),
    @("M2", "Col9a3", "Mag", "M2", 3, 1, 1.715756666666667, 5.14727, 0.3229633061525894, 0.3608812612892173, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 1.079239539277778, 9.7131558535, 0.1425984731593032, 0.1813726518443302),
    @("M2", "Col9a3", "Mag", "sCs", 3, 1, 1.715756666666667, 5.14727, 0.3229633061525894, 0.3608812612892173, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.8907659571200001, 5.34459574272, 0.1176957114753226, 0.09979902696030017),
    @("Neutro", "Col9a3", "Mag", "M1", 3, 1, 0.642518, 1.927554, 0.1209435705971609, 0.1351431183371138, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.17761767592, 1.59855908328, 0.02346838534960026, 0.02984971150852696),
    @("Neutro", "Col9a3", "Mag", "M2", 3, 1, 0.642518, 1.927554, 0.1209435705971609, 0.1351431183371138, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 0.4041545306333333, 3.6373907757, 0.0534003961968398, 0.06792058325153839),
    @("Neutro", "Col9a3", "Mag", "sCs", 3, 1, 0.642518, 1.927554, 0.1209435705971609, 0.1351431183371138, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.3335747850240001, 2.001448710144, 0.04407478905072085, 0.0373728235770485),
    @("sCs", "Col9a3", "Mag", "M1", 2, 1, 0.439315, 0.87863, 0.08269390852379505, 0.06160180107251903, 2, 0.6666666666666666, 0.27644, 0.8293199999999999, 0.1940440920813295, 0.2208748168298663, 0.1214442386, 0.7286654316, 0.01604626440015632, 0.0136062865282825),
    @("sCs", "Col9a3", "Mag", "M2", 2, 1, 0.439315, 0.87863, 0.08269390852379505, 0.06160180107251903, 3, 1, 0.6290166666666667, 1.88705, 0.4415315004607062, 0.502582625643659, 0.2763364569166667, 1.6580187415, 0.0365119655094716, 0.03095999492740498),
    @("sCs", "Col9a3", "Mag", "sCs", 2, 1, 0.439315, 0.87863, 0.08269390852379505, 0.06160180107251903, 2, 1, 0.5191680000000001, 1.038336, 0.3644244074579644, 0.2765425575264748, 0.22807828992, 0.9123131596800002, 0.03013567861416713, 0.01703551961683155),
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowValues = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Count; $j++) {
        $c = $j + 1
        $ws.Cells.Item($r, $c).Value = $rowValues[$j]
    }
}

